# Added exception handling to handle the hangs when the vpn isn't working well.
# Added a time out to the stream reader to further prevent hanging.
# (Data-side effect: two new "Named" method log rows were appended below the
# existing data, and column A was widened slightly to fit the new dates.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 3
$ws.Cells.Item(3, 1).Value = 42600.829097222224
$ws.Cells.Item(3, 2).Value = "Named"
$ws.Cells.Item(3, 3).Value = 7660
$ws.Cells.Item(3, 4).Value = 3671
$ws.Cells.Item(3, 5).Value = 193
$ws.Cells.Item(3, 6).Value = 49
$ws.Cells.Item(3, 7).Value = 37
$ws.Cells.Item(3, 8).Value = 56
$ws.Cells.Item(3, 9).Value = 42
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).Value = 0

# New row 4
$ws.Cells.Item(4, 1).Value = 42600.881828703707
$ws.Cells.Item(4, 2).Value = "Named"
$ws.Cells.Item(4, 3).Value = 7991
$ws.Cells.Item(4, 4).Value = 3681
$ws.Cells.Item(4, 5).Value = 194
$ws.Cells.Item(4, 6).Value = 49
$ws.Cells.Item(4, 7).Value = 37
$ws.Cells.Item(4, 8).Value = 56
$ws.Cells.Item(4, 9).Value = 42
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = 0

# Match the date formatting/style used by the existing date column (A2)
$ws.Range("A3:A4").NumberFormat = "m/d/yy h:mm"

# Column A was widened slightly to accommodate the new dates
$ws.Columns.Item(1).ColumnWidth = 14
